$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 22 (old summary row); this shifts the summary row
# (old 22 -> 25) and the footer row (old 23 -> 26) down automatically.
$ws.Rows("22:24").Insert()

# Copy formatting from the last 3 existing data rows (19:21) onto the newly inserted rows (22:24)
$ws.Range("A19:Q21").Copy()
$ws.Range("A22:Q24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-create the merged cells for the 3 new rows (same pattern used by every other data row)
$ws.Range("A22:B22").Merge()
$ws.Range("C22:G22").Merge()
$ws.Range("H22:K22").Merge()
$ws.Range("L22:M22").Merge()
$ws.Range("N22:O22").Merge()
$ws.Range("A23:B23").Merge()
$ws.Range("C23:G23").Merge()
$ws.Range("H23:K23").Merge()
$ws.Range("L23:M23").Merge()
$ws.Range("N23:O23").Merge()
$ws.Range("A24:B24").Merge()
$ws.Range("C24:G24").Merge()
$ws.Range("H24:K24").Merge()
$ws.Range("L24:M24").Merge()
$ws.Range("N24:O24").Merge()

# Row heights matching the regenerated report layout
$ws.Rows("22:22").RowHeight = 25.5
$ws.Rows("23:23").RowHeight = 24.75
$ws.Rows("24:24").RowHeight = 25.5
$ws.Rows("25:25").RowHeight = 24.75

# Write the full data table (rows 7-24) with the updated shortage list
# (3 new items were inserted in sorted order: DIAMONRECTA, TAMSULIN, and a new Arabic item)
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "B-COM I.M./I.V. 6 AMP"
$ws.Range("H7").Value = "2:3"
$ws.Range("L7").Value = "1"
$ws.Range("N7").Value = "48.00"
$ws.Range("P7").Value = "7.6800"
$ws.Range("Q7").Value = "0:1"
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "DANSET 4MG/2ML 3 AMP"
$ws.Range("H8").Value = "0:2"
$ws.Range("L8").Value = "1"
$ws.Range("N8").Value = "82.50"
$ws.Range("P8").Value = "27.2250"
$ws.Range("Q8").Value = "0:1"
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "DECLOPHEN 75MG/3ML 3 AMPOULES"
$ws.Range("H9").Value = "3:1"
$ws.Range("L9").Value = "1"
$ws.Range("N9").Value = "36.00"
$ws.Range("P9").Value = "11.8800"
$ws.Range("Q9").Value = "0:1"
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "DIAMONRECTA 5 MG 30 F.C. TAB."
$ws.Range("H10").Value = "1:1"
$ws.Range("L10").Value = "1"
$ws.Range("N10").Value = "187.50"
$ws.Range("P10").Value = "61.8750"
$ws.Range("Q10").Value = "0:1"
$ws.Range("A11").Value = 5
$ws.Range("C11").Value = "GASTROTIDINE 20 MG/2ML 3 I.M./I.V. AMPS"
$ws.Range("H11").Value = "2:2"
$ws.Range("L11").Value = "1"
$ws.Range("N11").Value = "54.00"
$ws.Range("P11").Value = "17.8200"
$ws.Range("Q11").Value = "0:1"
$ws.Range("A12").Value = 6
$ws.Range("C12").Value = "NANAZOXID 500MG 18 F.C. TABS."
$ws.Range("H12").Value = "0:2"
$ws.Range("L12").Value = "1"
$ws.Range("N12").Value = "114.00"
$ws.Range("P12").Value = "37.6200"
$ws.Range("Q12").Value = "0:1"
$ws.Range("A13").Value = 7
$ws.Range("C13").Value = "OSTEOCARE 30 TABS"
$ws.Range("H13").Value = "1:0"
$ws.Range("L13").Value = "1"
$ws.Range("N13").Value = "150.00"
$ws.Range("P13").Value = "75.0000"
$ws.Range("Q13").Value = "0:1"
$ws.Range("A14").Value = 8
$ws.Range("C14").Value = "TAMSULIN 0.4MG 28 CAPS"
$ws.Range("H14").Value = "1:1"
$ws.Range("L14").Value = "1"
$ws.Range("N14").Value = "124.00"
$ws.Range("P14").Value = "124.0000"
$ws.Range("Q14").Value = "1:0"
$ws.Range("A15").Value = 9
$ws.Range("C15").Value = "TOBRIN 0.3% EYE DROPS 5 ML"
$ws.Range("H15").Value = "2:0"
$ws.Range("L15").Value = "1"
$ws.Range("N15").Value = "23.00"
$ws.Range("P15").Value = "23.0000"
$ws.Range("Q15").Value = "1:0"
$ws.Range("A16").Value = 10
$ws.Range("C16").Value = "VASTAREL MR 35MG 30 F.C.TAB."
$ws.Range("H16").Value = "1:2"
$ws.Range("L16").Value = "1"
$ws.Range("N16").Value = "175.00"
$ws.Range("P16").Value = "175.0000"
$ws.Range("Q16").Value = "1:0"
$ws.Range("A17").Value = 11
$ws.Range("C17").Value = "VISCERALGINE 5MG/2ML IM IV 6 AMPOULES"
$ws.Range("H17").Value = "0:4"
$ws.Range("L17").Value = "1"
$ws.Range("N17").Value = "90.00"
$ws.Range("P17").Value = "14.4000"
$ws.Range("Q17").Value = "0:1"
$ws.Range("A18").Value = 12
$ws.Range("C18").Value = "WATER FOR INJECTION AMP. 5 ML"
$ws.Range("H18").Value = "8703:0"
$ws.Range("L18").Value = "1"
$ws.Range("N18").Value = "2.00"
$ws.Range("P18").Value = "2.0000"
$ws.Range("Q18").Value = "1:0"
$ws.Range("A19").Value = 13
$ws.Range("C19").Value = "بلاستر مترسيلك 2.5 سم"
$ws.Range("H19").Value = "27:0"
$ws.Range("L19").Value = "0"
$ws.Range("N19").Value = "25.00"
$ws.Range("P19").Value = "25.0000"
$ws.Range("Q19").Value = "1:0"
$ws.Range("A20").Value = 14
$ws.Range("C20").Value = "سرنجات 3 سم"
$ws.Range("H20").Value = "0:0"
$ws.Range("L20").Value = "0"
$ws.Range("N20").Value = "2.00"
$ws.Range("P20").Value = "6.0000"
$ws.Range("Q20").Value = "3:0"
$ws.Range("A21").Value = 15
$ws.Range("C21").Value = "سرنجات 5 سم"
$ws.Range("H21").Value = "0:0"
$ws.Range("L21").Value = "0"
$ws.Range("N21").Value = "3.00"
$ws.Range("P21").Value = "3.0000"
$ws.Range("Q21").Value = "1:0"
$ws.Range("A22").Value = 16
$ws.Range("C22").Value = "فرشه شعر اطفال الجو"
$ws.Range("H22").Value = "0:0"
$ws.Range("L22").Value = "0"
$ws.Range("N22").Value = "25.00"
$ws.Range("P22").Value = "25.0000"
$ws.Range("Q22").Value = "1:0"
$ws.Range("A23").Value = 17
$ws.Range("C23").Value = "قطن 100 جم"
$ws.Range("H23").Value = "23:0"
$ws.Range("L23").Value = "0"
$ws.Range("N23").Value = "20.00"
$ws.Range("P23").Value = "20.0000"
$ws.Range("Q23").Value = "1:0"
$ws.Range("A24").Value = 18
$ws.Range("C24").Value = "محلول ملح"
$ws.Range("H24").Value = "29:0"
$ws.Range("L24").Value = "0"
$ws.Range("N24").Value = "24.00"
$ws.Range("P24").Value = "24.0000"
$ws.Range("Q24").Value = "1:0"

# Update the totals row (sum of the "sale price" column) and the footer timestamp
$ws.Range("P25").Value = 680.5
$ws.Range("A26").Value = "Tuesday, 3 June, 2025 10:33 AM"
